# "Fruta / hortaliza, semanal" weekly update:
# Insert two new rows at the top of the data block (row 9) pushing the
# existing historical rows down by two, then populate the two new rows
# with this week's price data for "Alcachofa" at
# "Mapocho Venta Directa de Santiago".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 9 (shifts old rows 9-14 -> 11-16)
$ws.Range("A9:R10").EntireRow.Insert()

# --- New row 9: Primera quality, new weekly data ---
$ws.Cells.Item(9, 1).Value  = 12
$ws.Cells.Item(9, 2).Value  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(9, 3).Value  = "Metropolitana"
$ws.Cells.Item(9, 4).Value  = 44841
$ws.Cells.Item(9, 5).Value  = 13
$ws.Cells.Item(9, 6).Value  = 100112013
$ws.Cells.Item(9, 7).Value  = "Alcachofa"
$ws.Cells.Item(9, 8).Value  = "Española"
$ws.Cells.Item(9, 9).Value  = "Primera"
$ws.Cells.Item(9, 10).Value = 45
$ws.Cells.Item(9, 11).Value = 12000
$ws.Cells.Item(9, 12).Value = 12000
$ws.Cells.Item(9, 13).Value = 12000
$ws.Cells.Item(9, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 400
$ws.Cells.Item(9, 17).Value = 30
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# --- New row 10: Segunda quality, new weekly data ---
$ws.Cells.Item(10, 1).Value  = 12
$ws.Cells.Item(10, 2).Value  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(10, 3).Value  = "Metropolitana"
$ws.Cells.Item(10, 4).Value  = 44841
$ws.Cells.Item(10, 5).Value  = 13
$ws.Cells.Item(10, 6).Value  = 100112013
$ws.Cells.Item(10, 7).Value  = "Alcachofa"
$ws.Cells.Item(10, 8).Value  = "Española"
$ws.Cells.Item(10, 9).Value  = "Segunda"
$ws.Cells.Item(10, 10).Value = 45
$ws.Cells.Item(10, 11).Value = 10000
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 13).Value = 10000
$ws.Cells.Item(10, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 250
$ws.Cells.Item(10, 17).Value = 40
$ws.Cells.Item(10, 18).Value = "Hortaliza"
